$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column, copying the header style from the neighboring "sum" header (G1)
# so the new header cell H1 gets the same bold/centered/bordered formatting.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the data values for the new "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
